# Scheduled market-data refresh for Chocobo_Profits.xlsx.
# Columns H-N on each touched leve row hold scraped marketboard figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) - plain cached
# numbers, no formulas - so the runner just overwrites them with the newly
# fetched values. A handful of rows flip between "no profitable recipe"
# (blank M/N) and "has one" (populated M/N) as prices move, so some cells
# are cleared out and others newly populated alongside the refreshed ones.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 9833.333000000001
$ws.Cells.Item(113, 9).Value = 3750
$ws.Cells.Item(113, 11).Value = 3750
$ws.Cells.Item(113, 13).Value = -496

$ws.Cells.Item(129, 8).Value = 868.52576
$ws.Cells.Item(129, 10).Value = 884.2234
$ws.Cells.Item(129, 12).Value = 2652.6702
$ws.Cells.Item(129, 14).Value = -12652.6702

$ws.Cells.Item(137, 8).Value = 3178224.2
$ws.Cells.Item(137, 9).Value = 5954183
$ws.Cells.Item(137, 10).Value = 5700
$ws.Cells.Item(137, 11).Value = 17862549
$ws.Cells.Item(137, 12).Value = 17100
$ws.Cells.Item(137, 13).Value = -17859999
$ws.Cells.Item(137, 14).Value = -22200

$ws.Cells.Item(138, 8).Value = 2522.15
$ws.Cells.Item(138, 9).Value = 675.5625
$ws.Cells.Item(138, 10).Value = 2873.8809
$ws.Cells.Item(138, 11).Value = 2026.6875
$ws.Cells.Item(138, 12).Value = 8621.6427
$ws.Cells.Item(138, 13).Value = 3113.3125
$ws.Cells.Item(138, 14).Value = -18901.6427

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 17318440
$ws.Cells.Item(63, 9).Value = 46172504
$ws.Cells.Item(63, 11).Value = 46172504
$ws.Cells.Item(63, 13).Value = -46171818

$ws.Cells.Item(66, 8).Value = 17318440
$ws.Cells.Item(66, 9).Value = 46172504
$ws.Cells.Item(66, 11).Value = 230862520
$ws.Cells.Item(66, 13).Value = -230859088

$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = ""

$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = ""

$ws.Cells.Item(132, 8).Value = 3466.7273
$ws.Cells.Item(132, 9).Value = 1352
$ws.Cells.Item(132, 11).Value = 4056
$ws.Cells.Item(132, 13).Value = -1526

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 34532.25
$ws.Cells.Item(35, 10).Value = 34532.25
$ws.Cells.Item(35, 12).Value = 34532.25
$ws.Cells.Item(35, 14).Value = -35152.25

$ws.Cells.Item(62, 8).Value = 1500
$ws.Cells.Item(62, 9).Value = 1500
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 1500
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -814
$ws.Cells.Item(62, 14).Value = ""

$ws.Cells.Item(63, 8).Value = 29000
$ws.Cells.Item(63, 10).Value = 29000
$ws.Cells.Item(63, 12).Value = 29000
$ws.Cells.Item(63, 14).Value = -30372

$ws.Cells.Item(65, 8).Value = 1500
$ws.Cells.Item(65, 9).Value = 1500
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 4500
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -1068
$ws.Cells.Item(65, 14).Value = ""

$ws.Cells.Item(66, 8).Value = 29000
$ws.Cells.Item(66, 10).Value = 29000
$ws.Cells.Item(66, 12).Value = 87000
$ws.Cells.Item(66, 14).Value = -93864

$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).Value = ""

$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).Value = ""

$ws.Cells.Item(134, 8).Value = 1875.4166
$ws.Cells.Item(134, 9).Value = 1241.1666
$ws.Cells.Item(134, 10).Value = 5046.6665
$ws.Cells.Item(134, 11).Value = 3723.4998
$ws.Cells.Item(134, 12).Value = 15139.9995
$ws.Cells.Item(134, 13).Value = -1188.4998
$ws.Cells.Item(134, 14).Value = -20209.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2874
$ws.Cells.Item(31, 9).Value = 1022.95
$ws.Cells.Item(31, 11).Value = 1022.95
$ws.Cells.Item(31, 13).Value = -727.95

$ws.Cells.Item(34, 8).Value = 2874
$ws.Cells.Item(34, 9).Value = 1022.95
$ws.Cells.Item(34, 11).Value = 1022.95
$ws.Cells.Item(34, 13).Value = -820.95

$ws.Cells.Item(41, 8).Value = 38050.4
$ws.Cells.Item(41, 10).Value = 43813
$ws.Cells.Item(41, 12).Value = 43813
$ws.Cells.Item(41, 14).Value = -44669

$ws.Cells.Item(50, 8).Value = 31514.4
$ws.Cells.Item(50, 10).Value = 31514.4
$ws.Cells.Item(50, 12).Value = 31514.4
$ws.Cells.Item(50, 14).Value = -32764.4

$ws.Cells.Item(51, 8).Value = 32322.285
$ws.Cells.Item(51, 10).Value = 32322.285
$ws.Cells.Item(51, 12).Value = 32322.285
$ws.Cells.Item(51, 14).Value = -33794.285

$ws.Cells.Item(58, 8).Value = 3186.5
$ws.Cells.Item(58, 9).Value = 1804.6171
$ws.Cells.Item(58, 10).Value = 9090.909
$ws.Cells.Item(58, 11).Value = 1804.6171
$ws.Cells.Item(58, 12).Value = 9090.909
$ws.Cells.Item(58, 13).Value = -1601.6171
$ws.Cells.Item(58, 14).Value = -9496.909

$ws.Cells.Item(59, 8).Value = 31821.572
$ws.Cells.Item(59, 10).Value = 31821.572
$ws.Cells.Item(59, 12).Value = 31821.572
$ws.Cells.Item(59, 14).Value = -34111.572

$ws.Cells.Item(60, 8).Value = 14838.613
$ws.Cells.Item(60, 9).Value = 636.5
$ws.Cells.Item(60, 10).Value = 15818.069
$ws.Cells.Item(60, 11).Value = 636.5
$ws.Cells.Item(60, 12).Value = 15818.069
$ws.Cells.Item(60, 13).Value = -125.5
$ws.Cells.Item(60, 14).Value = -16840.069

$ws.Cells.Item(61, 8).Value = 32322.285
$ws.Cells.Item(61, 10).Value = 32322.285
$ws.Cells.Item(61, 12).Value = 32322.285
$ws.Cells.Item(61, 14).Value = -33018.285

$ws.Cells.Item(68, 8).Value = 57006
$ws.Cells.Item(68, 10).Value = 57006
$ws.Cells.Item(68, 12).Value = 57006
$ws.Cells.Item(68, 14).Value = -58504

$ws.Cells.Item(71, 8).Value = 57006
$ws.Cells.Item(71, 10).Value = 57006
$ws.Cells.Item(71, 12).Value = 171018
$ws.Cells.Item(71, 14).Value = -178506

$ws.Cells.Item(74, 8).Value = 32762.223
$ws.Cells.Item(74, 10).Value = 32762.223
$ws.Cells.Item(74, 12).Value = 32762.223
$ws.Cells.Item(74, 14).Value = -34510.223

$ws.Cells.Item(77, 8).Value = 32762.223
$ws.Cells.Item(77, 10).Value = 32762.223
$ws.Cells.Item(77, 12).Value = 98286.66900000001
$ws.Cells.Item(77, 14).Value = -107022.669

$ws.Cells.Item(132, 8).Value = 3543.8572
$ws.Cells.Item(132, 9).Value = 2068.4443
$ws.Cells.Item(132, 11).Value = 6205.3329
$ws.Cells.Item(132, 13).Value = -3675.3329

$ws.Cells.Item(134, 8).Value = 8611.875
$ws.Cells.Item(134, 9).Value = 11089.2
$ws.Cells.Item(134, 10).Value = 4483
$ws.Cells.Item(134, 11).Value = 33267.60000000001
$ws.Cells.Item(134, 12).Value = 13449
$ws.Cells.Item(134, 13).Value = -30732.60000000001
$ws.Cells.Item(134, 14).Value = -18519

$ws.Cells.Item(136, 8).Value = 3186.5
$ws.Cells.Item(136, 9).Value = 1804.6171
$ws.Cells.Item(136, 10).Value = 9090.909
$ws.Cells.Item(136, 11).Value = 5413.8513
$ws.Cells.Item(136, 12).Value = 27272.727
$ws.Cells.Item(136, 13).Value = -2863.8513
$ws.Cells.Item(136, 14).Value = -32372.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 534910.25
$ws.Cells.Item(5, 9).Value = 617.6667
$ws.Cells.Item(5, 11).Value = 1853.0001
$ws.Cells.Item(5, 13).Value = -1741.0001

$ws.Cells.Item(56, 8).Value = 4388.4614
$ws.Cells.Item(56, 9).Value = 4388.4614
$ws.Cells.Item(56, 11).Value = 4388.4614
$ws.Cells.Item(56, 13).Value = -3858.4614

$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 13).Value = ""

$ws.Cells.Item(122, 8).Value = 3097.4722
$ws.Cells.Item(122, 9).Value = 1199
$ws.Cells.Item(122, 10).Value = 3639.8928
$ws.Cells.Item(122, 11).Value = 10791
$ws.Cells.Item(122, 12).Value = 32759.0352
$ws.Cells.Item(122, 13).Value = -8341
$ws.Cells.Item(122, 14).Value = -37659.0352

$ws.Cells.Item(131, 8).Value = 701.0909
$ws.Cells.Item(131, 9).Value = 301.78946
$ws.Cells.Item(131, 10).Value = 795.925
$ws.Cells.Item(131, 11).Value = 905.3683800000001
$ws.Cells.Item(131, 12).Value = 2387.775
$ws.Cells.Item(131, 13).Value = 4134.63162
$ws.Cells.Item(131, 14).Value = -12467.775

$ws.Cells.Item(135, 8).Value = 534910.25
$ws.Cells.Item(135, 9).Value = 617.6667
$ws.Cells.Item(135, 11).Value = 5559.0003
$ws.Cells.Item(135, 13).Value = -3024.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 7407926
$ws.Cells.Item(107, 9).Value = 278.9
$ws.Cells.Item(107, 10).Value = 22223220
$ws.Cells.Item(107, 11).Value = 278.9
$ws.Cells.Item(107, 12).Value = 22223220
$ws.Cells.Item(107, 13).Value = 1641.1
$ws.Cells.Item(107, 14).Value = -22227060

$ws.Cells.Item(113, 8).Value = 1641.2727
$ws.Cells.Item(113, 9).Value = 1611.375
$ws.Cells.Item(113, 10).Value = 1721
$ws.Cells.Item(113, 11).Value = 1611.375
$ws.Cells.Item(113, 12).Value = 1721
$ws.Cells.Item(113, 13).Value = 558.625
$ws.Cells.Item(113, 14).Value = -6061

$ws.Cells.Item(133, 8).Value = 41230
$ws.Cells.Item(133, 10).Value = 41230
$ws.Cells.Item(133, 12).Value = 41230
$ws.Cells.Item(133, 14).Value = -51350

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4316.2354
$ws.Cells.Item(7, 9).Value = 3297.3333
$ws.Cells.Item(7, 10).Value = 5462.5
$ws.Cells.Item(7, 11).Value = 3297.3333
$ws.Cells.Item(7, 12).Value = 5462.5
$ws.Cells.Item(7, 13).Value = -3185.3333
$ws.Cells.Item(7, 14).Value = -5686.5

$ws.Cells.Item(18, 8).Value = 20000
$ws.Cells.Item(18, 10).Value = 20000
$ws.Cells.Item(18, 12).Value = 20000
$ws.Cells.Item(18, 14).Value = -20344

$ws.Cells.Item(40, 8).Value = 4191.591
$ws.Cells.Item(40, 9).Value = 2900.8333
$ws.Cells.Item(40, 11).Value = 2900.8333
$ws.Cells.Item(40, 13).Value = -2764.8333

$ws.Cells.Item(61, 8).Value = 1869.5834
$ws.Cells.Item(61, 9).Value = 1771.1111
$ws.Cells.Item(61, 11).Value = 1771.1111
$ws.Cells.Item(61, 13).Value = -1569.1111

$ws.Cells.Item(113, 8).Value = 1869.5834
$ws.Cells.Item(113, 9).Value = 1771.1111
$ws.Cells.Item(113, 11).Value = 1771.1111
$ws.Cells.Item(113, 13).Value = 398.8888999999999

$ws.Cells.Item(122, 8).Value = 4610.636
$ws.Cells.Item(122, 9).Value = 1651.3334
$ws.Cells.Item(122, 10).Value = 8161.8
$ws.Cells.Item(122, 11).Value = 4954.0002
$ws.Cells.Item(122, 12).Value = 24485.4
$ws.Cells.Item(122, 13).Value = -2504.0002
$ws.Cells.Item(122, 14).Value = -29385.4

$ws.Cells.Item(126, 8).Value = 4316.2354
$ws.Cells.Item(126, 9).Value = 3297.3333
$ws.Cells.Item(126, 10).Value = 5462.5
$ws.Cells.Item(126, 11).Value = 9891.999899999999
$ws.Cells.Item(126, 12).Value = 16387.5
$ws.Cells.Item(126, 13).Value = -7421.999899999999
$ws.Cells.Item(126, 14).Value = -21327.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1048.1666
$ws.Cells.Item(113, 9).Value = 940
$ws.Cells.Item(113, 10).Value = 1102.25
$ws.Cells.Item(113, 11).Value = 2820
$ws.Cells.Item(113, 12).Value = 3306.75
$ws.Cells.Item(113, 13).Value = -650
$ws.Cells.Item(113, 14).Value = -7646.75

$ws.Cells.Item(126, 8).Value = 889100.25
$ws.Cells.Item(126, 9).Value = 1700.75
$ws.Cells.Item(126, 11).Value = 5102.25
$ws.Cells.Item(126, 13).Value = -2632.25

